$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value while forcing it to be stored as text, so that
# numeric-looking strings (e.g. "0.75" or "642,530,686,576") are not
# auto-converted into numbers by Excel. The original cell style is restored
# afterwards so no stray number-format style is left on the cell.
function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

# ---------------------------------------------------------------------------
# Rows 7-11 undergo a cyclic shift: new row 7 = old row 8, new row 8 = old
# row 9, new row 9 = old row 10, new row 10 = old row 11, new row 11 = old
# row 7. Capture all the old values first (using Value2, which reads
# reliably in this runtime), then write them back in the new order.
# ---------------------------------------------------------------------------
$rows = 7..11
$colsText = @("A","D","E","F")
$colsCoordLike = @("I","J")
$colsNum  = @("G","H")

$oldValues = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $colsText) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    foreach ($c in $colsCoordLike) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    foreach ($c in $colsNum) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $oldValues[$r] = $rowData
}

# new row 7 gets old row 8, new row 8 gets old row 9, new row 9 gets old row 10,
# new row 10 gets old row 11, new row 11 gets old row 7
$mapping = @{7 = 8; 8 = 9; 9 = 10; 10 = 11; 11 = 7}

foreach ($r in $rows) {
    $srcRow = $mapping[$r]
    $src = $oldValues[$srcRow]
    foreach ($c in $colsText) {
        $ws.Range("$c$r").Value = $src[$c]
    }
    foreach ($c in $colsCoordLike) {
        Set-TextValue $ws.Range("$c$r") $src[$c]
    }
    foreach ($c in $colsNum) {
        $ws.Range("$c$r").Value = $src[$c]
    }
}

# ---------------------------------------------------------------------------
# Row 16 updates
# ---------------------------------------------------------------------------
$ws.Range("D16").Value = "image_20250807111344_ppp0.jpg"
Set-TextValue $ws.Range("I16") "642,530,686,576"
Set-TextValue $ws.Range("J16") "0.75"

# ---------------------------------------------------------------------------
# Row 17 updates
# ---------------------------------------------------------------------------
$ws.Range("D17").Value = "image_20250807111344_ppp0.jpg"
Set-TextValue $ws.Range("I17") "794,481,831,526"
